# Jogos_do_Dia_Betfair_Back_Lay_2025-12-29.xlsx update
# The feed refreshed: the "Italian Serie A" (Roma x Genoa) and
# "Portuguese Primeira Liga" (Porto x AVS Futebol SAD) fixtures dropped off
# the sheet, and the remaining "Friendly Matches" (Tlaxcala F.C x Pachuca)
# fixture moved up into row 2 with newly refreshed odds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old rows 3 (Portuguese Primeira Liga) and 4 (old Friendly
# Matches row) - row 2's data gets overwritten below with the refreshed
# Friendly Matches values, leaving just the header + one data row.
$ws.Rows("3:4").Delete()

# Refresh row 2 in place with the latest Friendly Matches odds.
$ws.Range("A2").Value = "Friendly Matches"
$ws.Range("B2").Value = "'2025-12-29"
$ws.Range("C2").Value = "21:00:00"
$ws.Range("D2").Value = "Tlaxcala F.C"
$ws.Range("E2").Value = "Pachuca"
$ws.Range("F2").Value = 4.4
$ws.Range("G2").Value = 5.4
$ws.Range("H2").Value = 1.71
$ws.Range("I2").Value = 1.79
$ws.Range("J2").Value = 4.3
$ws.Range("K2").Value = 5.1
$ws.Range("L2").Value = 1.31
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 4.7
$ws.Range("O2").Value = 1.23
$ws.Range("P2").Value = 2.32
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 1.6
$ws.Range("S2").Value = 2.6
$ws.Range("T2").Value = 1.64
$ws.Range("U2").Value = 2.28
$ws.Range("V2").Value = 2.24
$ws.Range("W2").Value = 1.23
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 19.5
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 990
$ws.Range("AC2").Value = 970
$ws.Range("AD2").Value = 12
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 23
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 12
